# Updated cryptos list on Sat Sep  9 22:35:32 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, [string]$Text)
    # Force the cell to keep a literal text value even when the
    # string looks numeric (e.g. "1.00"), without leaving a
    # lingering explicit number format/style on the cell.
    $Range.NumberFormat = "@"
    $Range.Value = $Text
    $Range.Style = "Normal"
}

$ws.Range("D2").Value = '26.023.58'
$ws.Range("E2").Value = '  +0.28%  '
$ws.Range("D3").Value = '1.643.91'
$ws.Range("E3").Value = '  +0.38%  '
$ws.Range("E4").Value = '  +0.27%  '
Set-TextValue $ws.Range("D5") '215.44'
$ws.Range("E5").Value = '  +0.29%  '
$ws.Range("E6").Value = '  +0.03%  '
Set-TextValue $ws.Range("D7") '1.00'
$ws.Range("E7").Value = '  +0.27%  '
$ws.Range("B8").Value = 'Dogecoin'
$ws.Range("C8").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
Set-TextValue $ws.Range("D8") '0.0639'
$ws.Range("E8").Value = '  +0.32%  '
$ws.Range("B9").Value = 'Cardano'
$ws.Range("C9").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
Set-TextValue $ws.Range("D9") '0.255'
$ws.Range("E9").Value = '  +0.15%  '
$ws.Range("E10").Value = '  -0.11%  '
$ws.Range("E11").Value = '  +0.29%  '
$ws.Range("E12").Value = '  +0.19%  '
$ws.Range("D13").Value = '1.591.17'
$ws.Range("E13").Value = '  -2.88%  '
Set-TextValue $ws.Range("D14") '0.543'
$ws.Range("E14").Value = '  -0.08%  '
Set-TextValue $ws.Range("D15") '63.43'
$ws.Range("E15").Value = '  +1.46%  '
$ws.Range("D16").Value = '0.0₃0760'
$ws.Range("E16").Value = '  +0.40%  '
$ws.Range("D17").Value = '26.059.28'
$ws.Range("E17").Value = '  +0.40%  '
$ws.Range("E18").Value = '  +0.29%  '
Set-TextValue $ws.Range("D19") '194.46'
$ws.Range("E19").Value = '  +0.38%  '
$ws.Range("E20").Value = '  -0.42%  '
Set-TextValue $ws.Range("D21") '9.90'
$ws.Range("E21").Value = '  -0.22%  '
$ws.Range("E22").Value = '  -1.00%  '
Set-TextValue $ws.Range("D23") '0.133'
$ws.Range("E23").Value = '  +4.50%  '
Set-TextValue $ws.Range("D24") '143.91'
$ws.Range("E24").Value = '  -0.21%  '
Set-TextValue $ws.Range("D26") '1.00'
$ws.Range("E26").Value = '  -0.04%  '
$ws.Range("E27").Value = '  +0.51%  '
$ws.Range("E29").Value = '  +0.30%  '
Set-TextValue $ws.Range("D30") '0.0495'
$ws.Range("E30").Value = '  -1.11%  '
$ws.Range("E31").Value = '  +1.09%  '
$ws.Range("E33").Value = '  -0.15%  '
$ws.Range("E34").Value = '  +1.00%  '
Set-TextValue $ws.Range("D35") '0.904'
$ws.Range("E35").Value = '  +0.31%  '
$ws.Range("D36").Value = '1.130.65'
Set-TextValue $ws.Range("D37") '0.539'
$ws.Range("E37").Value = '  -1.38%  '
$ws.Range("E38").Value = '  +0.20%  '
$ws.Range("E39").Value = '  -0.08%  '
Set-TextValue $ws.Range("D40") '5.44'
$ws.Range("E40").Value = '  +0.49%  '
Set-TextValue $ws.Range("D41") '98.87'
$ws.Range("E41").Value = '  -0.40%  '
$ws.Range("E42").Value = '  -0.56%  '
$ws.Range("E43").Value = '  +1.25%  '
Set-TextValue $ws.Range("D44") '56.47'
$ws.Range("E44").Value = '  +0.03%  '
$ws.Range("E45").Value = '  +2.52%  '
$ws.Range("E46").Value = '  -1.41%  '
Set-TextValue $ws.Range("D47") '7.77'
$ws.Range("E47").Value = '  +1.53%  '
$ws.Range("E48").Value = '  -0.25%  '
$ws.Range("E49").Value = '  +0.15%  '
$ws.Range("E50").Value = '  -1.33%  '
Set-TextValue $ws.Range("D51") '5.54'
$ws.Range("E51").Value = '  +0.09%  '
